$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.1477750351608889
$ws.Range("J2").Value = 0.1477750351608889
$ws.Range("M2").Value = 17.08155333333333
$ws.Range("N2").Value = 51.24466
$ws.Range("O2").Value = 0.3501540759902865
$ws.Range("P2").Value = 0.3501540759902865
$ws.Range("Q2").Value = 1.245569787513333
$ws.Range("R2").Value = 11.21012808762
$ws.Range("S2").Value = 0.05174403089119315
$ws.Range("T2").Value = 0.05174403089119315

# Row 3
$ws.Range("I3").Value = 0.1477750351608889
$ws.Range("J3").Value = 0.1477750351608889
$ws.Range("O3").Value = 0.2142771237573249
$ws.Range("P3").Value = 0.2142771237573249
$ws.Range("S3").Value = 0.03166480949741283
$ws.Range("T3").Value = 0.03166480949741283

# Row 4
$ws.Range("I4").Value = 0.1477750351608889
$ws.Range("J4").Value = 0.1477750351608889
$ws.Range("M4").Value = 8.398122666666666
$ws.Range("N4").Value = 25.194368
$ws.Range("O4").Value = 0.1721527793764119
$ws.Range("P4").Value = 0.1721527793764119
$ws.Range("Q4").Value = 0.6123827067306665
$ws.Range("R4").Value = 5.511444360575999
$ws.Range("S4").Value = 0.02543988302539402
$ws.Range("T4").Value = 0.02543988302539402

# Row 5
$ws.Range("I5").Value = 0.1477750351608889
$ws.Range("J5").Value = 0.1477750351608889
$ws.Range("M5").Value = 4.514486333333333
$ws.Range("N5").Value = 13.543459
$ws.Range("O5").Value = 0.09254227409953211
$ws.Range("P5").Value = 0.09254227409953213
$ws.Range("Q5").Value = 0.3291918289403333
$ws.Range("R5").Value = 2.962726460463
$ws.Range("S5").Value = 0.01367543780892698
$ws.Range("T5").Value = 0.01367543780892698

# Row 6
$ws.Range("I6").Value = 0.1477750351608889
$ws.Range("J6").Value = 0.1477750351608889
$ws.Range("M6").Value = 8.335727666666667
$ws.Range("N6").Value = 25.007183
$ws.Range("O6").Value = 0.1708737467764446
$ws.Range("P6").Value = 0.1708737467764446
$ws.Range("Q6").Value = 0.6078329257256667
$ws.Range("R6").Value = 5.470496331531001
$ws.Range("S6").Value = 0.02525087393796193
$ws.Range("T6").Value = 0.02525087393796193

# Row 7
$ws.Range("G7").Value = 0.420527
$ws.Range("H7").Value = 1.261581
$ws.Range("I7").Value = 0.852224964839111
$ws.Range("J7").Value = 0.852224964839111
$ws.Range("M7").Value = 17.08155333333333
$ws.Range("N7").Value = 51.24466
$ws.Range("O7").Value = 0.3501540759902865
$ws.Range("P7").Value = 0.3501540759902865
$ws.Range("Q7").Value = 7.183254378606667
$ws.Range("R7").Value = 64.64928940746
$ws.Range("S7").Value = 0.2984100450990933
$ws.Range("T7").Value = 0.2984100450990933

# Row 8
$ws.Range("G8").Value = 0.420527
$ws.Range("H8").Value = 1.261581
$ws.Range("I8").Value = 0.852224964839111
$ws.Range("J8").Value = 0.852224964839111
$ws.Range("O8").Value = 0.2142771237573249
$ws.Range("P8").Value = 0.2142771237573249
$ws.Range("Q8").Value = 4.395799429471
$ws.Range("R8").Value = 39.562194865239
$ws.Range("S8").Value = 0.182612314259912
$ws.Range("T8").Value = 0.182612314259912

# Row 9
$ws.Range("G9").Value = 0.420527
$ws.Range("H9").Value = 1.261581
$ws.Range("I9").Value = 0.852224964839111
$ws.Range("J9").Value = 0.852224964839111
$ws.Range("M9").Value = 8.398122666666666
$ws.Range("N9").Value = 25.194368
$ws.Range("O9").Value = 0.1721527793764119
$ws.Range("P9").Value = 0.1721527793764119
$ws.Range("Q9").Value = 3.531637330645333
$ws.Range("R9").Value = 31.784735975808
$ws.Range("S9").Value = 0.1467128963510179
$ws.Range("T9").Value = 0.1467128963510179

# Row 10
$ws.Range("G10").Value = 0.420527
$ws.Range("H10").Value = 1.261581
$ws.Range("I10").Value = 0.852224964839111
$ws.Range("J10").Value = 0.852224964839111
$ws.Range("M10").Value = 4.514486333333333
$ws.Range("N10").Value = 13.543459
$ws.Range("O10").Value = 0.09254227409953211
$ws.Range("P10").Value = 0.09254227409953213
$ws.Range("Q10").Value = 1.898463394297667
$ws.Range("R10").Value = 17.086170548679
$ws.Range("S10").Value = 0.07886683629060513
$ws.Range("T10").Value = 0.07886683629060515

# Row 11
$ws.Range("G11").Value = 0.420527
$ws.Range("H11").Value = 1.261581
$ws.Range("I11").Value = 0.852224964839111
$ws.Range("J11").Value = 0.852224964839111
$ws.Range("M11").Value = 8.335727666666667
$ws.Range("N11").Value = 25.007183
$ws.Range("O11").Value = 0.1708737467764446
$ws.Range("P11").Value = 0.1708737467764446
$ws.Range("Q11").Value = 3.505398548480334
$ws.Range("R11").Value = 31.548586936323
$ws.Range("S11").Value = 0.1456228728384827
$ws.Range("T11").Value = 0.1456228728384827
